$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) — rows keyed by their F column "想去人数" value
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1210
$wsExhibit.Range("F6").Value = 221
$wsExhibit.Range("F8").Value = 67
$wsExhibit.Range("F10").Value = 5682
$wsExhibit.Range("F11").Value = 5029
$wsExhibit.Range("F12").Value = 23
$wsExhibit.Range("F13").Value = 50
$wsExhibit.Range("F16").Value = 208
$wsExhibit.Range("F17").Value = 11

# Sheet "全部类型" (sheet4) — same underlying rows, plus extra rows, so F19 here
# corresponds to F17 on "展览"
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1210
$wsAll.Range("F6").Value = 221
$wsAll.Range("F8").Value = 67
$wsAll.Range("F10").Value = 5682
$wsAll.Range("F11").Value = 5029
$wsAll.Range("F12").Value = 23
$wsAll.Range("F13").Value = 50
$wsAll.Range("F16").Value = 208
$wsAll.Range("F19").Value = 11
